# Add a new reference row (row 5) to the tracking sheet:
#   A5 = paper title, B5 = authors, C5 = note on how it relates to
#   the F.Sener paper already tracked in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "An end-to-end generative framework for video segmentation and recognition`n"
$ws.Range("B5").Value = "Hilde Kuehne`nJuergen Gall`nThomas Serre`n"
$ws.Range("C5").Value = "Sener那篇文章用了这篇文章提供的特征"

# Title/author columns wrap, matching the style used by the row above (s="1").
$ws.Range("A5:B5").WrapText = $true

# Row sized to fit the wrapped text, as in the source row.
$ws.Rows.Item(5).RowHeight = 56

# Move the active selection to the newly added last cell, like Excel does
# right after data entry.
$ws.Range("C5").Select()
